$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.451.95'
$ws.Range("E2").Value = '  -3.68%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.952.15'
$ws.Range("E3").Value = '  -2.47%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.012'
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '321.58'
$ws.Range("E5").Value = '  -2.32%  '
$ws.Range("E6").Value = '  +0.17%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4776'
$ws.Range("E7").Value = '  -4.51%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4026'
$ws.Range("E8").Value = '  -4.36%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '53.58'
$ws.Range("E9").Value = '  -1.10%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08484'
$ws.Range("E10").Value = '  -5.65%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.059'
$ws.Range("E11").Value = '  -5.28%  '
$ws.Range("E12").Value = '  -4.97%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.944.61'
$ws.Range("E13").Value = '  -3.24%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.627'
$ws.Range("E14").Value = '  -5.12%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.192'
$ws.Range("E15").Value = '  -4.27%  '
$ws.Range("E16").Value = '  +0.17%  '
$ws.Range("E17").Value = '  -3.04%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '89.09'
$ws.Range("E18").Value = '  -5.51%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06604'
$ws.Range("E19").Value = '  -0.91%  '
$ws.Range("E20").Value = '  -4.61%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.013'
$ws.Range("E21").Value = '  +0.17%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.818'
$ws.Range("E22").Value = '  -2.39%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '28.474.60'
$ws.Range("E23").Value = '  -3.72%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.290'
$ws.Range("E25").Value = '  -0.39%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.182.29'
$ws.Range("E26").Value = '  -2.97%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '154.03'
$ws.Range("E27").Value = '  -3.09%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.17'
$ws.Range("E28").Value = '  -2.55%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.956'
$ws.Range("E29").Value = '  -6.89%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.158'
$ws.Range("E30").Value = '  -6.10%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '123.60'
$ws.Range("E31").Value = '  -3.52%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9948'
$ws.Range("E32").Value = '  -5.67%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09566'
$ws.Range("E33").Value = '  -3.93%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.447'
$ws.Range("E34").Value = '  -7.55%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.597'
$ws.Range("E35").Value = '  -4.01%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.661'
$ws.Range("E36").Value = '  -3.59%  '
$ws.Range("E37").Value = '  -5.07%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.789'
$ws.Range("E39").Value = '  -5.50%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.260'
$ws.Range("E40").Value = '  -3.58%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6241'
$ws.Range("E41").Value = '  -4.77%  '
$ws.Range("E42").Value = '  -4.86%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.012'
$ws.Range("E43").Value = '  +0.14%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.1929'
$ws.Range("E44").Value = '  -5.96%  '
$ws.Range("E45").Value = '  +1.88%  '
$ws.Range("E46").Value = '  -6.18%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '12.93'
$ws.Range("E47").Value = '  -4.26%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.065'
$ws.Range("E48").Value = '  -5.66%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.405'
$ws.Range("E49").Value = '  -2.95%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00000000331'
$ws.Range("E50").Value = '  -0.86%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06805'
$ws.Range("E51").Value = '  -2.66%  '
